# Remove the three decorative "Icon" pictures (Picture 14 / id 15, Picture 28 /
# id 29, Picture 5 / id 6) from the two "Demo" slides (slide 6 = sldId 357,
# slide 7 = sldId 392) - mirroring the same deletion already applied earlier
# to the other "Demo" slide (slide 1 / sldId 300).

$p = $ppt.ActivePresentation

$targetNames = @("Picture 14", "Picture 28", "Picture 5")

foreach ($slideIndex in 6, 7) {
    $s = $p.Slides.Item($slideIndex)
    for ($j = $s.Shapes.Count; $j -ge 1; $j--) {
        $sh = $s.Shapes.Item($j)
        if ($targetNames -contains $sh.Name) {
            $sh.Delete()
        }
    }
}
